# Daily attendance processing - 2025-11-22 19:42:26
#
# Normalizes the "Recorded By" column (G) on the active sheet: wherever the
# comma-separated list of recorders contains the literal entry "System" but
# it isn't already the first entry, move it to the front of the list while
# preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $v = $cell.Value2
    if ($v -eq $null) { continue }

    $parts = $v -split ", "
    if ($parts.Length -le 1) { continue }
    if ($parts[0] -eq "System") { continue }

    $sysIdx = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -eq "System") {
            $sysIdx = $i
        }
    }
    if ($sysIdx -lt 1) { continue }

    $newParts = @()
    $newParts += "System"
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $sysIdx) {
            $newParts += $parts[$i]
        }
    }

    $cell.Value = $newParts -join ", "
}
